# Apply scraped-schedule update for Línea 141 (horarios-141-2026-01-18.xlsx)
# New scrape at 13:14:31 adds rows and re-sorts each sheet by Hora_Llegada (col B).
$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:14:31"
$ws.Cells.Item(3, 1).Value = "Total filas: 176"
$ws.Cells.Item(38, 1).Value = "06:34:35"
$ws.Cells.Item(38, 2).Value = "08:29"
$ws.Cells.Item(38, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(38, 4).Value = 115
$ws.Cells.Item(38, 5).Value = "LP1912"
$ws.Cells.Item(39, 1).Value = "06:34:35"
$ws.Cells.Item(39, 2).Value = "08:29"
$ws.Cells.Item(39, 3).Value = "15_ABASTO"
$ws.Cells.Item(39, 4).Value = 115
$ws.Cells.Item(39, 5).Value = "LP1912"
$ws.Cells.Item(49, 1).Value = "08:52:40"
$ws.Cells.Item(49, 2).Value = "08:52"
$ws.Cells.Item(49, 3).Value = "215B_EL PATO"
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = "LP1912"
$ws.Cells.Item(50, 1).Value = "07:13:03"
$ws.Cells.Item(50, 2).Value = "08:52"
$ws.Cells.Item(50, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(50, 4).Value = 99
$ws.Cells.Item(50, 5).Value = "LP1912"
$ws.Cells.Item(137, 1).Value = "11:46:32"
$ws.Cells.Item(137, 2).Value = "12:34"
$ws.Cells.Item(137, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(137, 4).Value = 48
$ws.Cells.Item(137, 5).Value = "LP1912"
$ws.Cells.Item(138, 1).Value = "10:36:50"
$ws.Cells.Item(138, 2).Value = "12:34"
$ws.Cells.Item(138, 3).Value = "15_ABASTO"
$ws.Cells.Item(138, 4).Value = 118
$ws.Cells.Item(138, 5).Value = "LP1912"
$ws.Cells.Item(160, 1).Value = "12:33:02"
$ws.Cells.Item(160, 2).Value = "13:33"
$ws.Cells.Item(160, 3).Value = "14_ABASTO"
$ws.Cells.Item(160, 4).Value = 60
$ws.Cells.Item(160, 5).Value = "LP1912"
$ws.Cells.Item(161, 1).Value = "11:46:32"
$ws.Cells.Item(161, 2).Value = "13:33"
$ws.Cells.Item(161, 3).Value = "215A_EL PATO"
$ws.Cells.Item(161, 4).Value = 107
$ws.Cells.Item(161, 5).Value = "LP1912"
$ws.Cells.Item(166, 1).Value = "13:14:31"
$ws.Cells.Item(166, 2).Value = "14:02"
$ws.Cells.Item(166, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(166, 4).Value = 48
$ws.Cells.Item(166, 5).Value = "LP1912"
$ws.Cells.Item(167, 1).Value = "13:14:31"
$ws.Cells.Item(167, 2).Value = "14:05"
$ws.Cells.Item(167, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(167, 4).Value = 51
$ws.Cells.Item(167, 5).Value = "LP1912"
$ws.Cells.Item(168, 1).Value = "12:46:07"
$ws.Cells.Item(168, 2).Value = "14:08"
$ws.Cells.Item(168, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(168, 4).Value = 82
$ws.Cells.Item(168, 5).Value = "LP1912"
$ws.Cells.Item(169, 1).Value = "12:53:26"
$ws.Cells.Item(169, 2).Value = "14:09"
$ws.Cells.Item(169, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(169, 4).Value = 76
$ws.Cells.Item(169, 5).Value = "LP1912"
$ws.Cells.Item(170, 1).Value = "12:53:26"
$ws.Cells.Item(170, 2).Value = "14:16"
$ws.Cells.Item(170, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(170, 4).Value = 83
$ws.Cells.Item(170, 5).Value = "LP1912"
$ws.Cells.Item(171, 1).Value = "12:33:02"
$ws.Cells.Item(171, 2).Value = "14:17"
$ws.Cells.Item(171, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(171, 4).Value = 104
$ws.Cells.Item(171, 5).Value = "LP1912"
$ws.Cells.Item(172, 1).Value = "12:53:26"
$ws.Cells.Item(172, 2).Value = "14:17"
$ws.Cells.Item(172, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(172, 4).Value = 84
$ws.Cells.Item(172, 5).Value = "LP1912"
$ws.Cells.Item(173, 1).Value = "12:33:02"
$ws.Cells.Item(173, 2).Value = "14:18"
$ws.Cells.Item(173, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(173, 4).Value = 105
$ws.Cells.Item(173, 5).Value = "LP1912"
$ws.Cells.Item(174, 1).Value = "12:53:26"
$ws.Cells.Item(174, 2).Value = "14:27"
$ws.Cells.Item(174, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(174, 4).Value = 94
$ws.Cells.Item(174, 5).Value = "LP1912"
$ws.Cells.Item(175, 1).Value = "12:33:02"
$ws.Cells.Item(175, 2).Value = "14:32"
$ws.Cells.Item(175, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(175, 4).Value = 119
$ws.Cells.Item(175, 5).Value = "LP1912"
$ws.Cells.Item(176, 1).Value = "12:46:07"
$ws.Cells.Item(176, 2).Value = "14:34"
$ws.Cells.Item(176, 3).Value = "215C_EL PATO"
$ws.Cells.Item(176, 4).Value = 108
$ws.Cells.Item(176, 5).Value = "LP1912"
$ws.Cells.Item(177, 1).Value = "12:46:07"
$ws.Cells.Item(177, 2).Value = "14:39"
$ws.Cells.Item(177, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(177, 4).Value = 113
$ws.Cells.Item(177, 5).Value = "LP1912"
$ws.Cells.Item(178, 1).Value = "12:53:26"
$ws.Cells.Item(178, 2).Value = "14:47"
$ws.Cells.Item(178, 3).Value = "215B_EL PATO"
$ws.Cells.Item(178, 4).Value = 114
$ws.Cells.Item(178, 5).Value = "LP1912"
$ws.Cells.Item(179, 1).Value = "13:14:31"
$ws.Cells.Item(179, 2).Value = "14:54"
$ws.Cells.Item(179, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(179, 4).Value = 100
$ws.Cells.Item(179, 5).Value = "LP1912"
$ws.Cells.Item(180, 1).Value = "13:14:31"
$ws.Cells.Item(180, 2).Value = "15:02"
$ws.Cells.Item(180, 3).Value = "10_OLMOS"
$ws.Cells.Item(180, 4).Value = 108
$ws.Cells.Item(180, 5).Value = "LP1912"
$ws.Cells.Item(181, 1).Value = "13:14:31"
$ws.Cells.Item(181, 2).Value = "15:13"
$ws.Cells.Item(181, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(181, 4).Value = 119
$ws.Cells.Item(181, 5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:14:31"
$ws.Cells.Item(3, 1).Value = "Total filas: 31"
$ws.Cells.Item(36, 1).Value = "13:14:31"
$ws.Cells.Item(36, 2).Value = "14:54"
$ws.Cells.Item(36, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(36, 4).Value = 100
$ws.Cells.Item(36, 5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:14:31"
$ws.Cells.Item(3, 1).Value = "Total filas: 27"
$ws.Cells.Item(27, 1).Value = "13:14:31"
$ws.Cells.Item(27, 2).Value = "13:16"
$ws.Cells.Item(27, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = "L6203"
$ws.Cells.Item(28, 1).Value = "11:33:52"
$ws.Cells.Item(28, 2).Value = "13:20"
$ws.Cells.Item(28, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(28, 4).Value = 107
$ws.Cells.Item(28, 5).Value = "L6173"
$ws.Cells.Item(29, 1).Value = "11:46:32"
$ws.Cells.Item(29, 2).Value = "13:21"
$ws.Cells.Item(29, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(29, 4).Value = 95
$ws.Cells.Item(29, 5).Value = "L6173"
$ws.Cells.Item(30, 1).Value = "12:11:21"
$ws.Cells.Item(30, 2).Value = "13:57"
$ws.Cells.Item(30, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(30, 4).Value = 106
$ws.Cells.Item(30, 5).Value = "L6203"
$ws.Cells.Item(31, 1).Value = "13:14:31"
$ws.Cells.Item(31, 2).Value = "14:03"
$ws.Cells.Item(31, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(31, 4).Value = 49
$ws.Cells.Item(31, 5).Value = "L6203"
$ws.Cells.Item(32, 1).Value = "12:46:07"
$ws.Cells.Item(32, 2).Value = "14:27"
$ws.Cells.Item(32, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(32, 4).Value = 101
$ws.Cells.Item(32, 5).Value = "L6203"
